$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value()
    if ($val -ne $null) {
        $newVal = $val -replace '\(F\)$', ''
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
